$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "MEC-1A-Gestão"
$ws.Range("C3").Value = "MEC-2A-Mecanica material"
$ws.Range("E3").Value = "-"

# Row 4
$ws.Range("C4").Value = "MEC-2A-Mecanica material"
$ws.Range("D4").Value = "-"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("D6").Value = "MCT-3A-Usinagem"
$ws.Range("F6").Value = "MEC-3A-Usinagem"

# Row 7
$ws.Range("D7").Value = "MCT-3A-Usinagem"
$ws.Range("F7").Value = "MEC-3A-Usinagem"

$wb.Save()
